$wb = $excel.ActiveWorkbook

# The f9878666-90b1-48cc-b4fa-4a6a87ee1180 file finished handback, so its
# status moves from "Ready for handoff" / the stale-version error message
# to "Handed back: in sync with en-US", with fresh handback timestamps and
# a cleared error detail, on all three sheets (Overview, zh-cn, de-de).

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the f9878666 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 5).Value = $handedBack   # E3 (zh-cn status column)
$wsOverview.Cells.Item(3, 6).Value = $handedBack   # F3 (de-de status column)

# --- zh-cn sheet: row 3 is the f9878666 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value = $handedBack              # C3 Status
$wsZhCn.Cells.Item(3, 11).Value = "2016-08-25 22:47:48"   # K3 Latest Handback DateTime
$wsZhCn.Cells.Item(3, 16).Value = ""                      # P3 Error Detail

# --- de-de sheet: row 3 is the f9878666 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value = $handedBack              # C3 Status
$wsDeDe.Cells.Item(3, 11).Value = "2016-08-25 22:47:54"   # K3 Latest Handback DateTime
$wsDeDe.Cells.Item(3, 16).Value = ""                      # P3 Error Detail
